$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 98-99 (existing rows 98:115 shift down to 100:117),
# carrying the style (e.g. date format on column D) of the old row 98 down
# with them, matching Excel's normal "insert rows" shift behavior.
$ws.Rows("98:99").Insert()

# Fill in the data for the two newly inserted rows (new Mandarina / Murcott
# price entries dated 2021-09-10, Provincia de Limari origin).
$ws.Range("A98").Value = 7
$ws.Range("B98").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C98").Value = "Ñuble"
$ws.Range("D98").Value = 44449
$ws.Range("E98").Value = 16
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100102
$ws.Range("H98").Value = "Cítricos"
$ws.Range("I98").Value = 100102004
$ws.Range("J98").Value = "Mandarina"
$ws.Range("K98").Value = "Murcott"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 240
$ws.Range("N98").Value = 6500
$ws.Range("O98").Value = 7000
$ws.Range("P98").Value = 6750
$ws.Range("Q98").Value = "`$/bandeja 10 kilos"
$ws.Range("R98").Value = "Provincia de Limarí"
$ws.Range("S98").Value = 675
$ws.Range("T98").Value = 10

$ws.Range("A99").Value = 7
$ws.Range("B99").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C99").Value = "Ñuble"
$ws.Range("D99").Value = 44449
$ws.Range("E99").Value = 16
$ws.Range("F99").Value = "Fruta"
$ws.Range("G99").Value = 100102
$ws.Range("H99").Value = "Cítricos"
$ws.Range("I99").Value = 100102004
$ws.Range("J99").Value = "Mandarina"
$ws.Range("K99").Value = "Murcott"
$ws.Range("L99").Value = "Segunda"
$ws.Range("M99").Value = 90
$ws.Range("N99").Value = 6000
$ws.Range("O99").Value = 6000
$ws.Range("P99").Value = 6000
$ws.Range("Q99").Value = "`$/bandeja 10 kilos"
$ws.Range("R99").Value = "Provincia de Limarí"
$ws.Range("S99").Value = 600
$ws.Range("T99").Value = 10
